# Add a new "Sampled Residues" column (F) with values for each data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Sampled Residues"
$ws.Range("F2").Value = "APR:FED:SK:KED"
$ws.Range("F4").Value = "WT"
$ws.Range("F3").Value = "IS::SRT:VN"

# Select the new full data range and zoom in, matching the saved view state.
$ws.Range("A1:F4").Select()
$excel.ActiveWindow.Zoom = 186
